$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "274.93"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-1.19%"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "26.76"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-2.50%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.892"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "2.19%"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06330"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1.24%"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.862"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-0.84%"

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.315"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "1.50%"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.254"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "33.33%"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8693"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-1.21%"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1709"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "17.90%"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.05031"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-5.56%"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07419"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "1.03%"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.02961"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-4.79%"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09022"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.42%"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001574"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1.41%"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0006318"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.57%"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.005782"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-2.01%"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.447"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.05%"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1335"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "1.76%"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.911"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "1.64%"

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "1.06%"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001178"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.08%"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004249"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-0.82%"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001201"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.05%"

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-0.18%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04064"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "0.55%"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006731"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-1.64%"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1167"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "1.00%"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002181"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-0.41%"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01075"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-10.53%"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005291"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "3.93%"

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-29.66%"

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-37.48%"
